$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2:E2').NumberFormat = '@'
$ws.Range('D2').Value = '67.560.65'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D2:E2').ClearFormats()

$ws.Range('D3:E3').NumberFormat = '@'
$ws.Range('D3').Value = '2.636.16'
$ws.Range('E3').Value = '  -1.03%  '
$ws.Range('D3:E3').ClearFormats()

$ws.Range('D4:E4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.34%  '
$ws.Range('D4:E4').ClearFormats()

$ws.Range('D5:E5').NumberFormat = '@'
$ws.Range('D5').Value = '597.41'
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D5:E5').ClearFormats()

$ws.Range('D6:E6').NumberFormat = '@'
$ws.Range('D6').Value = '169.65'
$ws.Range('E6').Value = '  +3.55%  '
$ws.Range('D6:E6').ClearFormats()

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('E7').ClearFormats()

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.90%  '
$ws.Range('E8').ClearFormats()

$ws.Range('D9:E9').NumberFormat = '@'
$ws.Range('D9').Value = '2.631.18'
$ws.Range('E9').Value = '  -1.25%  '
$ws.Range('D9:E9').ClearFormats()

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.35%  '
$ws.Range('E10').ClearFormats()

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.46%  '
$ws.Range('E11').ClearFormats()

$ws.Range('D12:E12').NumberFormat = '@'
$ws.Range('D12').Value = '0.364'
$ws.Range('E12').Value = '  +2.27%  '
$ws.Range('D12:E12').ClearFormats()

$ws.Range('D13:E13').NumberFormat = '@'
$ws.Range('D13').Value = '5.25'
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('D13:E13').ClearFormats()

$ws.Range('D14:E14').NumberFormat = '@'
$ws.Range('D14').Value = '27.82'
$ws.Range('E14').Value = '  +0.54%  '
$ws.Range('D14:E14').ClearFormats()

$ws.Range('D15:E15').NumberFormat = '@'
$ws.Range('D15').Value = '3.099.76'
$ws.Range('E15').Value = '  -1.62%  '
$ws.Range('D15:E15').ClearFormats()

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000183'
$ws.Range('D16').ClearFormats()

$ws.Range('D17:E17').NumberFormat = '@'
$ws.Range('D17').Value = '66.961.80'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D17:E17').ClearFormats()

$ws.Range('D18:E18').NumberFormat = '@'
$ws.Range('D18').Value = '2.609.04'
$ws.Range('E18').Value = '  -1.94%  '
$ws.Range('D18:E18').ClearFormats()

$ws.Range('D19:E19').NumberFormat = '@'
$ws.Range('D19').Value = '12.19'
$ws.Range('E19').Value = '  +5.18%  '
$ws.Range('D19:E19').ClearFormats()

$ws.Range('D20:E20').NumberFormat = '@'
$ws.Range('D20').Value = '8.14'
$ws.Range('E20').Value = '  +8.80%  '
$ws.Range('D20:E20').ClearFormats()

$ws.Range('D21:E21').NumberFormat = '@'
$ws.Range('D21').Value = '358.04'
$ws.Range('E21').Value = '  -0.69%  '
$ws.Range('D21:E21').ClearFormats()

$ws.Range('D22:E22').NumberFormat = '@'
$ws.Range('D22').Value = '4.34'
$ws.Range('E22').Value = '  -0.59%  '
$ws.Range('D22:E22').ClearFormats()

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.77%  '
$ws.Range('E23').ClearFormats()

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +6.33%  '
$ws.Range('E24').ClearFormats()

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('E25').ClearFormats()

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.50%  '
$ws.Range('E26').ClearFormats()

$ws.Range('D27:E27').NumberFormat = '@'
$ws.Range('D27').Value = '69.86'
$ws.Range('E27').Value = '  -1.75%  '
$ws.Range('D27:E27').ClearFormats()

$ws.Range('D28:E28').NumberFormat = '@'
$ws.Range('D28').Value = '2.761.18'
$ws.Range('E28').Value = '  -1.41%  '
$ws.Range('D28:E28').ClearFormats()

$ws.Range('D29:E29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D29:E29').ClearFormats()

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.57%  '
$ws.Range('E30').ClearFormats()

$ws.Range('D31:E31').NumberFormat = '@'
$ws.Range('D31').Value = '551.96'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D31:E31').ClearFormats()

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('E32').ClearFormats()

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.32%  '
$ws.Range('E33').ClearFormats()

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('E34').ClearFormats()

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +6.70%  '
$ws.Range('E35').ClearFormats()

$ws.Range('D36:E36').NumberFormat = '@'
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D36:E36').ClearFormats()

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.64%  '
$ws.Range('E37').ClearFormats()

$ws.Range('D38:E38').NumberFormat = '@'
$ws.Range('D38').Value = '158.08'
$ws.Range('E38').Value = '  +2.50%  '
$ws.Range('D38:E38').ClearFormats()

$ws.Range('D39:E39').NumberFormat = '@'
$ws.Range('D39').Value = '19.06'
$ws.Range('E39').Value = '  -1.69%  '
$ws.Range('D39:E39').ClearFormats()

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.10%  '
$ws.Range('E40').ClearFormats()

$ws.Range('D41:E41').NumberFormat = '@'
$ws.Range('D41').Value = '5.21'
$ws.Range('E41').Value = '  -0.93%  '
$ws.Range('D41:E41').ClearFormats()

$ws.Range('B42:E42').NumberFormat = '@'
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').Value = '18.14'
$ws.Range('E42').Value = '  +1.41%  '
$ws.Range('B42:E42').ClearFormats()

$ws.Range('B43:E43').NumberFormat = '@'
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '1.80'
$ws.Range('E43').Value = '  -0.75%  '
$ws.Range('B43:E43').ClearFormats()

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('E44').ClearFormats()

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.79%  '
$ws.Range('E45').ClearFormats()

$ws.Range('D46:E46').NumberFormat = '@'
$ws.Range('D46').Value = '40.16'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('D46:E46').ClearFormats()

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.39%  '
$ws.Range('E47').ClearFormats()

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('E48').ClearFormats()

$ws.Range('D49:E49').NumberFormat = '@'
$ws.Range('D49').Value = '152.35'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D49:E49').ClearFormats()

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('E50').ClearFormats()

$ws.Range('D51:E51').NumberFormat = '@'
$ws.Range('D51').Value = '1.72'
$ws.Range('E51').Value = '  -0.44%  '
$ws.Range('D51:E51').ClearFormats()
